$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DispatchDate cell (B2) previously held a raw numeric date serial (367)
# formatted with a date number format. Change it to hold the literal date
# text "31-12-2018" instead (keeping the existing DD/MM/YYYY cell style/format).
$ws.Range("B2").Value = "31-12-2018"

# Update the active selection to B3 (as in the target workbook).
$ws.Range("B3").Select()
